# Apply weekly consolidated update to the Arándano (blue) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing date number format for column D (Fecha) before the
# bulk value assignment so new/reused date cells keep style index 2.
$ws.Range("D2:D32").NumberFormat = $ws.Range("D2").NumberFormat

$data = New-Object 'object[,]' 31,20
$data[0,0] = 4
$data[0,1] = "Feria Lagunitas de Puerto Montt"
$data[0,2] = "Los Lagos"
$data[0,3] = (Get-Date -Year 2022 -Month 11 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[0,4] = 10
$data[0,5] = "Fruta"
$data[0,6] = 100101
$data[0,7] = "Berries"
$data[0,8] = 100101001
$data[0,9] = "Arándano (blue)"
$data[0,10] = "Sin especificar"
$data[0,11] = "Primera"
$data[0,12] = 400
$data[0,13] = 7500
$data[0,14] = 8000
$data[0,15] = 7750
$data[0,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[0,17] = "Provincia de Curicó"
$data[0,18] = 5167
$data[0,19] = 1.5
$data[1,0] = 4
$data[1,1] = "Feria Lagunitas de Puerto Montt"
$data[1,2] = "Los Lagos"
$data[1,3] = (Get-Date -Year 2022 -Month 12 -Day 9 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[1,4] = 10
$data[1,5] = "Fruta"
$data[1,6] = 100101
$data[1,7] = "Berries"
$data[1,8] = 100101001
$data[1,9] = "Arándano (blue)"
$data[1,10] = "Sin especificar"
$data[1,11] = "Primera"
$data[1,12] = 400
$data[1,13] = 5000
$data[1,14] = 5500
$data[1,15] = 5250
$data[1,16] = "`$/bandeja 2 kilos"
$data[1,17] = "Provincia de Curicó"
$data[1,18] = 2625
$data[1,19] = 2
$data[2,0] = 4
$data[2,1] = "Feria Lagunitas de Puerto Montt"
$data[2,2] = "Los Lagos"
$data[2,3] = (Get-Date -Year 2021 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[2,4] = 10
$data[2,5] = "Fruta"
$data[2,6] = 100101
$data[2,7] = "Berries"
$data[2,8] = 100101001
$data[2,9] = "Arándano (blue)"
$data[2,10] = "Sin especificar"
$data[2,11] = "Primera"
$data[2,12] = 400
$data[2,13] = 5000
$data[2,14] = 5500
$data[2,15] = 5250
$data[2,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[2,17] = "Región del Maule"
$data[2,18] = 3500
$data[2,19] = 1.5
$data[3,0] = 4
$data[3,1] = "Feria Lagunitas de Puerto Montt"
$data[3,2] = "Los Lagos"
$data[3,3] = (Get-Date -Year 2022 -Month 12 -Day 5 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[3,4] = 10
$data[3,5] = "Fruta"
$data[3,6] = 100101
$data[3,7] = "Berries"
$data[3,8] = 100101001
$data[3,9] = "Arándano (blue)"
$data[3,10] = "Sin especificar"
$data[3,11] = "Primera"
$data[3,12] = 400
$data[3,13] = 5500
$data[3,14] = 6000
$data[3,15] = 5750
$data[3,16] = "`$/bandeja 2 kilos"
$data[3,17] = "Provincia de Curicó"
$data[3,18] = 2875
$data[3,19] = 2
$data[4,0] = 4
$data[4,1] = "Feria Lagunitas de Puerto Montt"
$data[4,2] = "Los Lagos"
$data[4,3] = (Get-Date -Year 2021 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[4,4] = 10
$data[4,5] = "Fruta"
$data[4,6] = 100101
$data[4,7] = "Berries"
$data[4,8] = 100101001
$data[4,9] = "Arándano (blue)"
$data[4,10] = "Sin especificar"
$data[4,11] = "Primera"
$data[4,12] = 300
$data[4,13] = 3700
$data[4,14] = 3800
$data[4,15] = 3750
$data[4,16] = "`$/kilo"
$data[4,17] = "Región del Maule"
$data[4,18] = 3750
$data[4,19] = 1
$data[5,0] = 4
$data[5,1] = "Feria Lagunitas de Puerto Montt"
$data[5,2] = "Los Lagos"
$data[5,3] = (Get-Date -Year 2021 -Month 11 -Day 19 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[5,4] = 10
$data[5,5] = "Fruta"
$data[5,6] = 100101
$data[5,7] = "Berries"
$data[5,8] = 100101001
$data[5,9] = "Arándano (blue)"
$data[5,10] = "Sin especificar"
$data[5,11] = "Primera"
$data[5,12] = 200
$data[5,13] = 3700
$data[5,14] = 3800
$data[5,15] = 3750
$data[5,16] = "`$/kilo"
$data[5,17] = "Región del Maule"
$data[5,18] = 3750
$data[5,19] = 1
$data[6,0] = 4
$data[6,1] = "Feria Lagunitas de Puerto Montt"
$data[6,2] = "Los Lagos"
$data[6,3] = (Get-Date -Year 2021 -Month 12 -Day 17 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[6,4] = 10
$data[6,5] = "Fruta"
$data[6,6] = 100101
$data[6,7] = "Berries"
$data[6,8] = 100101001
$data[6,9] = "Arándano (blue)"
$data[6,10] = "Sin especificar"
$data[6,11] = "Primera"
$data[6,12] = 400
$data[6,13] = 5000
$data[6,14] = 5500
$data[6,15] = 5250
$data[6,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[6,17] = "Región del Maule"
$data[6,18] = 3500
$data[6,19] = 1.5
$data[7,0] = 4
$data[7,1] = "Feria Lagunitas de Puerto Montt"
$data[7,2] = "Los Lagos"
$data[7,3] = (Get-Date -Year 2021 -Month 12 -Day 7 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[7,4] = 10
$data[7,5] = "Fruta"
$data[7,6] = 100101
$data[7,7] = "Berries"
$data[7,8] = 100101001
$data[7,9] = "Arándano (blue)"
$data[7,10] = "Sin especificar"
$data[7,11] = "Primera"
$data[7,12] = 400
$data[7,13] = 5000
$data[7,14] = 5500
$data[7,15] = 5250
$data[7,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[7,17] = "Región del Maule"
$data[7,18] = 3500
$data[7,19] = 1.5
$data[8,0] = 4
$data[8,1] = "Feria Lagunitas de Puerto Montt"
$data[8,2] = "Los Lagos"
$data[8,3] = (Get-Date -Year 2020 -Month 12 -Day 11 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[8,4] = 10
$data[8,5] = "Fruta"
$data[8,6] = 100101
$data[8,7] = "Berries"
$data[8,8] = 100101001
$data[8,9] = "Arándano (blue)"
$data[8,10] = "Sin especificar"
$data[8,11] = "Primera"
$data[8,12] = 300
$data[8,13] = 5000
$data[8,14] = 6000
$data[8,15] = 5500
$data[8,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[8,17] = "Provincia de Curicó"
$data[8,18] = 3667
$data[8,19] = 1.5
$data[9,0] = 4
$data[9,1] = "Feria Lagunitas de Puerto Montt"
$data[9,2] = "Los Lagos"
$data[9,3] = (Get-Date -Year 2021 -Month 12 -Day 21 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[9,4] = 10
$data[9,5] = "Fruta"
$data[9,6] = 100101
$data[9,7] = "Berries"
$data[9,8] = 100101001
$data[9,9] = "Arándano (blue)"
$data[9,10] = "Sin especificar"
$data[9,11] = "Primera"
$data[9,12] = 400
$data[9,13] = 5000
$data[9,14] = 5500
$data[9,15] = 5250
$data[9,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[9,17] = "Región del Maule"
$data[9,18] = 3500
$data[9,19] = 1.5
$data[10,0] = 4
$data[10,1] = "Feria Lagunitas de Puerto Montt"
$data[10,2] = "Los Lagos"
$data[10,3] = (Get-Date -Year 2022 -Month 11 -Day 11 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[10,4] = 10
$data[10,5] = "Fruta"
$data[10,6] = 100101
$data[10,7] = "Berries"
$data[10,8] = 100101001
$data[10,9] = "Arándano (blue)"
$data[10,10] = "Sin especificar"
$data[10,11] = "Primera"
$data[10,12] = 200
$data[10,13] = 7500
$data[10,14] = 8000
$data[10,15] = 7750
$data[10,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[10,17] = "Provincia de Curicó"
$data[10,18] = 5167
$data[10,19] = 1.5
$data[11,0] = 4
$data[11,1] = "Feria Lagunitas de Puerto Montt"
$data[11,2] = "Los Lagos"
$data[11,3] = (Get-Date -Year 2021 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[11,4] = 10
$data[11,5] = "Fruta"
$data[11,6] = 100101
$data[11,7] = "Berries"
$data[11,8] = 100101001
$data[11,9] = "Arándano (blue)"
$data[11,10] = "Sin especificar"
$data[11,11] = "Primera"
$data[11,12] = 400
$data[11,13] = 5000
$data[11,14] = 5500
$data[11,15] = 5250
$data[11,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[11,17] = "Región del Maule"
$data[11,18] = 3500
$data[11,19] = 1.5
$data[12,0] = 4
$data[12,1] = "Feria Lagunitas de Puerto Montt"
$data[12,2] = "Los Lagos"
$data[12,3] = (Get-Date -Year 2020 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[12,4] = 10
$data[12,5] = "Fruta"
$data[12,6] = 100101
$data[12,7] = "Berries"
$data[12,8] = 100101001
$data[12,9] = "Arándano (blue)"
$data[12,10] = "Sin especificar"
$data[12,11] = "Primera"
$data[12,12] = 400
$data[12,13] = 5500
$data[12,14] = 6000
$data[12,15] = 5750
$data[12,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[12,17] = "Provincia de Curicó"
$data[12,18] = 3833
$data[12,19] = 1.5
$data[13,0] = 4
$data[13,1] = "Feria Lagunitas de Puerto Montt"
$data[13,2] = "Los Lagos"
$data[13,3] = (Get-Date -Year 2022 -Month 12 -Day 6 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[13,4] = 10
$data[13,5] = "Fruta"
$data[13,6] = 100101
$data[13,7] = "Berries"
$data[13,8] = 100101001
$data[13,9] = "Arándano (blue)"
$data[13,10] = "Sin especificar"
$data[13,11] = "Primera"
$data[13,12] = 400
$data[13,13] = 5000
$data[13,14] = 5500
$data[13,15] = 5250
$data[13,16] = "`$/bandeja 2 kilos"
$data[13,17] = "Provincia de Curicó"
$data[13,18] = 2625
$data[13,19] = 2
$data[14,0] = 4
$data[14,1] = "Feria Lagunitas de Puerto Montt"
$data[14,2] = "Los Lagos"
$data[14,3] = (Get-Date -Year 2020 -Month 11 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[14,4] = 10
$data[14,5] = "Fruta"
$data[14,6] = 100101
$data[14,7] = "Berries"
$data[14,8] = 100101001
$data[14,9] = "Arándano (blue)"
$data[14,10] = "Sin especificar"
$data[14,11] = "Segunda"
$data[14,12] = 200
$data[14,13] = 6500
$data[14,14] = 7000
$data[14,15] = 6750
$data[14,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[14,17] = "Provincia de Curicó"
$data[14,18] = 4500
$data[14,19] = 1.5
$data[15,0] = 4
$data[15,1] = "Feria Lagunitas de Puerto Montt"
$data[15,2] = "Los Lagos"
$data[15,3] = (Get-Date -Year 2021 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[15,4] = 10
$data[15,5] = "Fruta"
$data[15,6] = 100101
$data[15,7] = "Berries"
$data[15,8] = 100101001
$data[15,9] = "Arándano (blue)"
$data[15,10] = "Sin especificar"
$data[15,11] = "Primera"
$data[15,12] = 160
$data[15,13] = 3600
$data[15,14] = 3700
$data[15,15] = 3650
$data[15,16] = "`$/kilo"
$data[15,17] = "Región del Maule"
$data[15,18] = 3650
$data[15,19] = 1
$data[16,0] = 4
$data[16,1] = "Feria Lagunitas de Puerto Montt"
$data[16,2] = "Los Lagos"
$data[16,3] = (Get-Date -Year 2020 -Month 12 -Day 1 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[16,4] = 10
$data[16,5] = "Fruta"
$data[16,6] = 100101
$data[16,7] = "Berries"
$data[16,8] = 100101001
$data[16,9] = "Arándano (blue)"
$data[16,10] = "Sin especificar"
$data[16,11] = "Primera"
$data[16,12] = 200
$data[16,13] = 6000
$data[16,14] = 6500
$data[16,15] = 6250
$data[16,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[16,17] = "Provincia de Curicó"
$data[16,18] = 4167
$data[16,19] = 1.5
$data[17,0] = 4
$data[17,1] = "Feria Lagunitas de Puerto Montt"
$data[17,2] = "Los Lagos"
$data[17,3] = (Get-Date -Year 2022 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[17,4] = 10
$data[17,5] = "Fruta"
$data[17,6] = 100101
$data[17,7] = "Berries"
$data[17,8] = 100101001
$data[17,9] = "Arándano (blue)"
$data[17,10] = "Sin especificar"
$data[17,11] = "Primera"
$data[17,12] = 400
$data[17,13] = 7500
$data[17,14] = 8000
$data[17,15] = 7750
$data[17,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[17,17] = "Provincia de Curicó"
$data[17,18] = 5167
$data[17,19] = 1.5
$data[18,0] = 4
$data[18,1] = "Feria Lagunitas de Puerto Montt"
$data[18,2] = "Los Lagos"
$data[18,3] = (Get-Date -Year 2022 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[18,4] = 10
$data[18,5] = "Fruta"
$data[18,6] = 100101
$data[18,7] = "Berries"
$data[18,8] = 100101001
$data[18,9] = "Arándano (blue)"
$data[18,10] = "Sin especificar"
$data[18,11] = "Primera"
$data[18,12] = 300
$data[18,13] = 7500
$data[18,14] = 8000
$data[18,15] = 7750
$data[18,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[18,17] = "Provincia de Curicó"
$data[18,18] = 5167
$data[18,19] = 1.5
$data[19,0] = 4
$data[19,1] = "Feria Lagunitas de Puerto Montt"
$data[19,2] = "Los Lagos"
$data[19,3] = (Get-Date -Year 2022 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[19,4] = 10
$data[19,5] = "Fruta"
$data[19,6] = 100101
$data[19,7] = "Berries"
$data[19,8] = 100101001
$data[19,9] = "Arándano (blue)"
$data[19,10] = "Sin especificar"
$data[19,11] = "Primera"
$data[19,12] = 400
$data[19,13] = 6000
$data[19,14] = 6500
$data[19,15] = 6250
$data[19,16] = "`$/bandeja 2 kilos"
$data[19,17] = "Provincia de Curicó"
$data[19,18] = 3125
$data[19,19] = 2
$data[20,0] = 4
$data[20,1] = "Feria Lagunitas de Puerto Montt"
$data[20,2] = "Los Lagos"
$data[20,3] = (Get-Date -Year 2022 -Month 11 -Day 29 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[20,4] = 10
$data[20,5] = "Fruta"
$data[20,6] = 100101
$data[20,7] = "Berries"
$data[20,8] = 100101001
$data[20,9] = "Arándano (blue)"
$data[20,10] = "Sin especificar"
$data[20,11] = "Primera"
$data[20,12] = 300
$data[20,13] = 8000
$data[20,14] = 8500
$data[20,15] = 8250
$data[20,16] = "`$/bandeja 2 kilos"
$data[20,17] = "Provincia de Curicó"
$data[20,18] = 4125
$data[20,19] = 2
$data[21,0] = 4
$data[21,1] = "Feria Lagunitas de Puerto Montt"
$data[21,2] = "Los Lagos"
$data[21,3] = (Get-Date -Year 2022 -Month 12 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[21,4] = 10
$data[21,5] = "Fruta"
$data[21,6] = 100101
$data[21,7] = "Berries"
$data[21,8] = 100101001
$data[21,9] = "Arándano (blue)"
$data[21,10] = "Sin especificar"
$data[21,11] = "Primera"
$data[21,12] = 400
$data[21,13] = 4000
$data[21,14] = 4400
$data[21,15] = 4200
$data[21,16] = "`$/bandeja 2 kilos"
$data[21,17] = "Provincia de Curicó"
$data[21,18] = 2100
$data[21,19] = 2
$data[22,0] = 4
$data[22,1] = "Feria Lagunitas de Puerto Montt"
$data[22,2] = "Los Lagos"
$data[22,3] = (Get-Date -Year 2022 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[22,4] = 10
$data[22,5] = "Fruta"
$data[22,6] = 100101
$data[22,7] = "Berries"
$data[22,8] = 100101001
$data[22,9] = "Arándano (blue)"
$data[22,10] = "Sin especificar"
$data[22,11] = "Primera"
$data[22,12] = 400
$data[22,13] = 4500
$data[22,14] = 5000
$data[22,15] = 4750
$data[22,16] = "`$/bandeja 2 kilos"
$data[22,17] = "Provincia de Curicó"
$data[22,18] = 2375
$data[22,19] = 2
$data[23,0] = 4
$data[23,1] = "Feria Lagunitas de Puerto Montt"
$data[23,2] = "Los Lagos"
$data[23,3] = (Get-Date -Year 2022 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[23,4] = 10
$data[23,5] = "Fruta"
$data[23,6] = 100101
$data[23,7] = "Berries"
$data[23,8] = 100101001
$data[23,9] = "Arándano (blue)"
$data[23,10] = "Sin especificar"
$data[23,11] = "Primera"
$data[23,12] = 200
$data[23,13] = 7500
$data[23,14] = 8000
$data[23,15] = 7750
$data[23,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[23,17] = "Provincia de Curicó"
$data[23,18] = 5167
$data[23,19] = 1.5
$data[24,0] = 4
$data[24,1] = "Feria Lagunitas de Puerto Montt"
$data[24,2] = "Los Lagos"
$data[24,3] = (Get-Date -Year 2022 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[24,4] = 10
$data[24,5] = "Fruta"
$data[24,6] = 100101
$data[24,7] = "Berries"
$data[24,8] = 100101001
$data[24,9] = "Arándano (blue)"
$data[24,10] = "Sin especificar"
$data[24,11] = "Primera"
$data[24,12] = 120
$data[24,13] = 8000
$data[24,14] = 8500
$data[24,15] = 8250
$data[24,16] = "`$/bandeja 2 kilos"
$data[24,17] = "Provincia de Curicó"
$data[24,18] = 4125
$data[24,19] = 2
$data[25,0] = 4
$data[25,1] = "Feria Lagunitas de Puerto Montt"
$data[25,2] = "Los Lagos"
$data[25,3] = (Get-Date -Year 2022 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[25,4] = 10
$data[25,5] = "Fruta"
$data[25,6] = 100101
$data[25,7] = "Berries"
$data[25,8] = 100101001
$data[25,9] = "Arándano (blue)"
$data[25,10] = "Sin especificar"
$data[25,11] = "Primera"
$data[25,12] = 300
$data[25,13] = 8000
$data[25,14] = 8500
$data[25,15] = 8250
$data[25,16] = "`$/bandeja 2 kilos"
$data[25,17] = "Provincia de Curicó"
$data[25,18] = 4125
$data[25,19] = 2
$data[26,0] = 4
$data[26,1] = "Feria Lagunitas de Puerto Montt"
$data[26,2] = "Los Lagos"
$data[26,3] = (Get-Date -Year 2021 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[26,4] = 10
$data[26,5] = "Fruta"
$data[26,6] = 100101
$data[26,7] = "Berries"
$data[26,8] = 100101001
$data[26,9] = "Arándano (blue)"
$data[26,10] = "Sin especificar"
$data[26,11] = "Primera"
$data[26,12] = 80
$data[26,13] = 3700
$data[26,14] = 3800
$data[26,15] = 3750
$data[26,16] = "`$/kilo"
$data[26,17] = "Región del Maule"
$data[26,18] = 3750
$data[26,19] = 1
$data[27,0] = 4
$data[27,1] = "Feria Lagunitas de Puerto Montt"
$data[27,2] = "Los Lagos"
$data[27,3] = (Get-Date -Year 2020 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[27,4] = 10
$data[27,5] = "Fruta"
$data[27,6] = 100101
$data[27,7] = "Berries"
$data[27,8] = 100101001
$data[27,9] = "Arándano (blue)"
$data[27,10] = "Sin especificar"
$data[27,11] = "Primera"
$data[27,12] = 100
$data[27,13] = 7000
$data[27,14] = 7000
$data[27,15] = 7000
$data[27,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[27,17] = "Provincia de Curicó"
$data[27,18] = 4667
$data[27,19] = 1.5
$data[28,0] = 4
$data[28,1] = "Feria Lagunitas de Puerto Montt"
$data[28,2] = "Los Lagos"
$data[28,3] = (Get-Date -Year 2020 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[28,4] = 10
$data[28,5] = "Fruta"
$data[28,6] = 100101
$data[28,7] = "Berries"
$data[28,8] = 100101001
$data[28,9] = "Arándano (blue)"
$data[28,10] = "Sin especificar"
$data[28,11] = "Segunda"
$data[28,12] = 100
$data[28,13] = 6500
$data[28,14] = 6500
$data[28,15] = 6500
$data[28,16] = "`$/bandeja 12 canastillos 125 gramos"
$data[28,17] = "Provincia de Curicó"
$data[28,18] = 4333
$data[28,19] = 1.5
$data[29,0] = 4
$data[29,1] = "Feria Lagunitas de Puerto Montt"
$data[29,2] = "Los Lagos"
$data[29,3] = (Get-Date -Year 2021 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[29,4] = 10
$data[29,5] = "Fruta"
$data[29,6] = 100101
$data[29,7] = "Berries"
$data[29,8] = 100101001
$data[29,9] = "Arándano (blue)"
$data[29,10] = "Sin especificar"
$data[29,11] = "Primera"
$data[29,12] = 400
$data[29,13] = 3500
$data[29,14] = 3600
$data[29,15] = 3550
$data[29,16] = "`$/kilo"
$data[29,17] = "Región del Maule"
$data[29,18] = 3550
$data[29,19] = 1
$data[30,0] = 4
$data[30,1] = "Feria Lagunitas de Puerto Montt"
$data[30,2] = "Los Lagos"
$data[30,3] = (Get-Date -Year 2022 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$data[30,4] = 10
$data[30,5] = "Fruta"
$data[30,6] = 100101
$data[30,7] = "Berries"
$data[30,8] = 100101001
$data[30,9] = "Arándano (blue)"
$data[30,10] = "Sin especificar"
$data[30,11] = "Primera"
$data[30,12] = 300
$data[30,13] = 8500
$data[30,14] = 9000
$data[30,15] = 8750
$data[30,16] = "`$/bandeja 2 kilos"
$data[30,17] = "Provincia de Curicó"
$data[30,18] = 4375
$data[30,19] = 2

$ws.Range("A2:T32").Value = $data
